$d = $word.ActiveDocument

# Target bullet currently reads:
#   "Tablet devices: Apple iPad 32GB, Wi-Fi and cellular connectivity "
# and needs to become:
#   "Tablet devices: Apple iPad 128GB, Wi-Fi and cellular connectivity "
# with the text split across three runs:
#   "Tablet devices: Apple iPad " | "128" | "GB, Wi-Fi and cellular connectivity "
#
# A plain Find/Replace (or a single InsertAfter/Text assignment) collapses
# the whole paragraph back into one run, so instead we perform the edit as
# two separate tracked-and-individually-accepted mutations (insert "128",
# then delete the old "32"). Accepting each revision right away keeps the
# surrounding untouched text split into its own runs instead of folding
# everything into a single run.

# --- Step 1: insert "128" immediately before "32GB" ---
$d.TrackRevisions = $true
$rng = $d.Content
$rng.Find.Execute("32GB", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$insertPos = $d.Range($rng.Start, $rng.Start)
$insertPos.InsertAfter("128")
$d.TrackRevisions = $false
$d.Revisions.Item(1).Accept()

# --- Step 2: delete the old "32" that now sits right before "GB" ---
$d.TrackRevisions = $true
$rng2 = $d.Content
$rng2.Find.Execute("12832GB", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$delRange = $d.Range($rng2.Start + 3, $rng2.Start + 5)
$delRange.Delete()
$d.TrackRevisions = $false
$d.Revisions.Item(1).Accept()
